# Scheduled-runner update: refresh currentAveragePrice / Leve profit
# calculations (columns H-N) on the per-craft profit sheets, per the
# latest market data pull. Only numeric value cells are touched --
# no structural / formatting changes.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1999.5
$ws.Range("I18").Value = 1999.5
$ws.Range("K18").Value = 1999.5
$ws.Range("M18").Value = -1715.5

$ws.Range("H28").Value = 773.38464
$ws.Range("I28").Value = 619.9091
$ws.Range("K28").Value = 619.9091
$ws.Range("M28").Value = -134.9091

$ws.Range("H33").Value = 586.7778
$ws.Range("I33").Value = 285.125
$ws.Range("K33").Value = 285.125
$ws.Range("M33").Value = -56.125

$ws.Range("H53").Value = 340.52
$ws.Range("I53").Value = 462.7857
$ws.Range("J53").Value = 184.90909
$ws.Range("K53").Value = 462.7857
$ws.Range("L53").Value = 184.90909
$ws.Range("M53").Value = 174.2143
$ws.Range("N53").Value = -1458.90909

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H74").Value = 65718.75
$ws.Range("I74").Value = 4150
$ws.Range("J74").Value = 168333.33
$ws.Range("K74").Value = 4150
$ws.Range("L74").Value = 168333.33
$ws.Range("M74").Value = -3214
$ws.Range("N74").Value = -170205.33

$ws.Range("H77").Value = 65718.75
$ws.Range("I77").Value = 4150
$ws.Range("J77").Value = 168333.33
$ws.Range("K77").Value = 20750
$ws.Range("L77").Value = 841666.6499999999
$ws.Range("M77").Value = -16070
$ws.Range("N77").Value = -851026.6499999999

$ws.Range("H107").Value = 1475.2972
$ws.Range("I107").Value = 640.7778
$ws.Range("K107").Value = 640.7778
$ws.Range("M107").Value = 1279.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2087.1052
$ws.Range("I2").Value = 1403
$ws.Range("K2").Value = 1403
$ws.Range("M2").Value = -1290

$ws.Range("H32").Value = 8188.4565
$ws.Range("I32").Value = 8188.4565
$ws.Range("K32").Value = 8188.4565
$ws.Range("M32").Value = -7901.4565

$ws.Range("H74").Value = 1782.4667
$ws.Range("I74").Value = 1195.7142
$ws.Range("K74").Value = 1195.7142
$ws.Range("M74").Value = -321.7141999999999

$ws.Range("H77").Value = 1782.4667
$ws.Range("I77").Value = 1195.7142
$ws.Range("K77").Value = 5978.571
$ws.Range("M77").Value = -1610.571

$ws.Range("H82").Value = 40181
$ws.Range("J82").Value = 40181
$ws.Range("L82").Value = 40181
$ws.Range("N82").Value = -40903

$ws.Range("H85").Value = 40181
$ws.Range("J85").Value = 40181
$ws.Range("L85").Value = 40181
$ws.Range("N85").Value = -42677

$ws.Range("H97").Value = 1082.0625
$ws.Range("I97").Value = 648
$ws.Range("K97").Value = 648
$ws.Range("M97").Value = -152

$ws.Range("H116").Value = 2087.1052
$ws.Range("I116").Value = 1403
$ws.Range("K116").Value = 1403
$ws.Range("M116").Value = 891

$ws.Range("H122").Value = 1998.6666
$ws.Range("I122").Value = 2278.4
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 6835.200000000001
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = -4385.200000000001
$ws.Range("N122").Value = -6700

$ws.Range("H132").Value = 3649.6667
$ws.Range("I132").Value = 3626.9092
$ws.Range("K132").Value = 10880.7276
$ws.Range("M132").Value = -8350.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2087.1052
$ws.Range("I3").Value = 1403
$ws.Range("K3").Value = 1403
$ws.Range("M3").Value = -1289

$ws.Range("H20").Value = 6800
$ws.Range("I20").Value = 6800
$ws.Range("K20").Value = 6800
$ws.Range("M20").Value = -6553

$ws.Range("H82").Value = 13038.667

$ws.Range("H85").Value = 13038.667

$ws.Range("H99").Value = 4555
$ws.Range("I99").Value = 4555
$ws.Range("K99").Value = 4555
$ws.Range("M99").Value = -3057

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 665.2727
$ws.Range("J2").Value = 766.3333
$ws.Range("L2").Value = 766.3333
$ws.Range("N2").Value = -992.3333

$ws.Range("H16").Value = 6834.8335
$ws.Range("I16").Value = 4201.8
$ws.Range("K16").Value = 4201.8
$ws.Range("M16").Value = -3914.8

$ws.Range("H22").Value = 893.625
$ws.Range("J22").Value = 894.8333
$ws.Range("L22").Value = 894.8333
$ws.Range("N22").Value = -1594.8333

$ws.Range("H58").Value = 2123.3
$ws.Range("I58").Value = 2407.5
$ws.Range("J58").Value = 986.5
$ws.Range("K58").Value = 2407.5
$ws.Range("L58").Value = 986.5
$ws.Range("M58").Value = -2204.5
$ws.Range("N58").Value = -1392.5

$ws.Range("H62").Value = 9500.4375
$ws.Range("I62").Value = 10410.8
$ws.Range("J62").Value = 7983.1665
$ws.Range("K62").Value = 10410.8
$ws.Range("L62").Value = 7983.1665
$ws.Range("M62").Value = -9786.799999999999
$ws.Range("N62").Value = -9231.166499999999

$ws.Range("H65").Value = 9500.4375
$ws.Range("I65").Value = 10410.8
$ws.Range("J65").Value = 7983.1665
$ws.Range("K65").Value = 52054
$ws.Range("L65").Value = 39915.8325
$ws.Range("M65").Value = -48934
$ws.Range("N65").Value = -46155.8325

$ws.Range("H113").Value = 6834.8335
$ws.Range("I113").Value = 4201.8
$ws.Range("K113").Value = 4201.8
$ws.Range("M113").Value = -2031.8

$ws.Range("H134").Value = 3805.6667
$ws.Range("I134").Value = 3853
$ws.Range("K134").Value = 11559
$ws.Range("M134").Value = -9024

$ws.Range("H136").Value = 2123.3
$ws.Range("I136").Value = 2407.5
$ws.Range("J136").Value = 986.5
$ws.Range("K136").Value = 7222.5
$ws.Range("L136").Value = 2959.5
$ws.Range("M136").Value = -4672.5
$ws.Range("N136").Value = -8059.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1630.2858
$ws.Range("I34").Value = 282.8
$ws.Range("J34").Value = 4999
$ws.Range("K34").Value = 848.4000000000001
$ws.Range("L34").Value = 14997
$ws.Range("M34").Value = -764.4000000000001
$ws.Range("N34").Value = -15165

$ws.Range("H55").Value = 224.5
$ws.Range("I55").Value = 224.5
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 673.5
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -496.5
$ws.Range("N55").ClearContents()

$ws.Range("H68").Value = 1220.5
$ws.Range("I68").Value = 960.6667
$ws.Range("K68").Value = 2882.0001
$ws.Range("M68").Value = -2071.0001

$ws.Range("H71").Value = 1220.5
$ws.Range("I71").Value = 960.6667
$ws.Range("K71").Value = 8646.0003
$ws.Range("M71").Value = -4590.0003

$ws.Range("H97").Value = 890.46155
$ws.Range("I97").Value = 1183.3334
$ws.Range("K97").Value = 3550.0002
$ws.Range("M97").Value = -3054.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 991
$ws.Range("I113").Value = 991
$ws.Range("K113").Value = 991
$ws.Range("M113").Value = 1179

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5744.3335
$ws.Range("I16").Value = 5744.3335
$ws.Range("K16").Value = 5744.3335
$ws.Range("M16").Value = -5574.3335

$ws.Range("H122").Value = 3002.5
$ws.Range("J122").Value = 4505
$ws.Range("L122").Value = 13515
$ws.Range("N122").Value = -18415

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 15000
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15298

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()

$ws.Range("H113").Value = 2084.4285
$ws.Range("I113").Value = 1016.5455
$ws.Range("K113").Value = 3049.6365
$ws.Range("M113").Value = -879.6364999999996

$ws.Range("H122").Value = 3040.4285
$ws.Range("I122").Value = 2821
$ws.Range("K122").Value = 8463
$ws.Range("M122").Value = -6013

$ws.Range("H126").Value = 3498.6
$ws.Range("I126").Value = 3498.6
$ws.Range("K126").Value = 10495.8
$ws.Range("M126").Value = -8025.799999999999

$ws.Range("H136").Value = 2237.8948
$ws.Range("I136").Value = 2140
$ws.Range("K136").Value = 6420
$ws.Range("M136").Value = -3870
